$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 0. I7 already carries the workbook's only "custom" style (the old date
#    format, s="1"). Re-pointing IT to the new hh:mm:ss format FIRST lets
#    the engine mutate that slot in place instead of retiring it and
#    allocating a new one later (which would leave an orphaned, unused xf
#    entry behind in styles.xml). I3:I6 then just copy the resulting format.
# ---------------------------------------------------------------------------
$ws.Range("I7").NumberFormat = "hh:mm:ss"
$ws.Range("I7").Copy()
$ws.Range("I3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("I5").PasteSpecial(-4122)
$ws.Range("I6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 1. Column data rewrite (A2:B7) - replace the old 8-row ramp with the new
#    6-row repeating sequence, and clear the now-unused rows 8:9 plus the
#    stray helper formulas that used to live in G8/I8/I10.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = 30
$ws.Range("B2").Value = 2.7
$ws.Range("A3").Value = 30
$ws.Range("B3").Value = 2.8
$ws.Range("A4").Value = 30
$ws.Range("B4").Value = 2.9
$ws.Range("A5").Value = 30
$ws.Range("B5").Value = 2.7
$ws.Range("A6").Value = 30
$ws.Range("B6").Value = 2.8
$ws.Range("A7").Value = 30
$ws.Range("B7").Value = 2.9

$ws.Range("A8:I9").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("I10").ClearContents()

# ---------------------------------------------------------------------------
# 2. Re-label the header row + add the new "Elapsed Runtime" column header,
#    and the new D-column labels / E-column values.
#    Writing D7 / B1 first just re-touches already-existing shared strings
#    (NUM STEPS / Speed (RPM)) without disturbing their slot; "Hour #" falls
#    out of the table entirely once A1 is overwritten below.
# ---------------------------------------------------------------------------
$ws.Range("D7").Value = "NUM STEPS"
$ws.Range("B1").Value = "Speed (RPM)"
$ws.Range("H1").Value = "Elapsed Runtime (seconds)"
$ws.Range("D2").Value = "Min Flow Rate (RPM"
$ws.Range("D3").Value = "Max Flow Rate (RPM)"
$ws.Range("D4").Value = "Offset Voltage (mV)"
$ws.Range("D10").Value = "RAW_STEPS"

# E2 / E3 are free-form text ("2.0" / "5.0"), right aligned, text numfmt.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").HorizontalAlignment = -4152
$ws.Range("E2").Value = "2.0"

# E3 picks up the same style by copy/paste (avoids leaving stray
# intermediate style entries behind in styles.xml); format must land BEFORE
# the text value, otherwise a still-General cell coerces "5.0" to 5.
$ws.Range("E2").Copy()
$ws.Range("E3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E3").Value = "5.0"

# E4 stays numeric (400) but shares the same text-numfmt/right-align style -
# set the value first (while still General), THEN paste the formatting, so
# it is NOT coerced into a text string by the "@" format.
$ws.Range("E4").Value = 400
$ws.Range("E2").Copy()
$ws.Range("E4").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A1").Value = "Step Duration (seconds)"
$ws.Range("D12").Value = "Sequence (copy this -->)"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Bold styling for the D-column labels (D2,D3,D4,D7,D10,D12) - set once on
#    D2, then fan out via copy/paste so no orphan style entries accumulate.
# ---------------------------------------------------------------------------
$ws.Range("D2").Font.Bold = $true
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("D12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row (A1,B1,H1) - bold, right aligned, wrapped.
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").HorizontalAlignment = -4152
$ws.Range("A1").WrapText = $true
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Running-total "Elapsed Runtime" column (H).
# ---------------------------------------------------------------------------
$ws.Range("H2").Value = 0
$ws.Range("H3").Formula = "=A3+H2"
$ws.Range("H4").Formula = "=A4+H3"
$ws.Range("H5").Formula = "=A5+H4"
$ws.Range("H6").Formula = "=A6+H5"
$ws.Range("H7").Formula = "=A7+H6"

# ---------------------------------------------------------------------------
# 5. Summary formulas.
# ---------------------------------------------------------------------------
$ws.Range("E7").Formula = "=COUNT(A2:A10000)"
$ws.Range("E10").Formula = '=TEXTJOIN(",","TRUE",A2:B10000)'
$ws.Range("E12").Formula = '=CONCAT("<",E2,",",E3,",",E4,",",E7,",",E10,">")'

# ---------------------------------------------------------------------------
# 6. Column widths / row height cosmetics.
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 29.25

$ws.Columns.Item(1).ColumnWidth = 15.71
$ws.Columns.Item(2).ColumnWidth = 12
$ws.Columns.Item(4).ColumnWidth = 25.57
$ws.Columns.Item(5).ColumnWidth = 11.43
$ws.Columns.Item(6).ColumnWidth = 35.57
$ws.Columns.Item(8).ColumnWidth = 16.86
$ws.Columns.Item(9).ColumnWidth = 23.29

# ---------------------------------------------------------------------------
# 7. Selection, matching the saved sheetView state.
# ---------------------------------------------------------------------------
$ws.Range("D28").Select()
